$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 13.417
$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.945
$ws.Range("A12").Value = -21.694
$ws.Range("C13").Value = -12.686
$ws.Range("A18").Value = -21.694
$ws.Range("E20").Value = 12.932
